# "Generate Report for Archive"
#
# The localization status report is regenerated:
#   - Items that were previously marked "Ready for handoff" are now
#     shown as "In Translation" (Overview sheet's per-locale status
#     columns, plus the Status column on each locale sheet).
#   - The "Status"-ish columns that used to be sized to fit the longer
#     "Ready for handoff" text are narrowed to fit "In Translation".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newText = "In Translation"

# --- Update every cell that shows the old status text -----------------
# Overview: columns E (zh-cn) and F (de-de), rows 2-4
foreach ($row in 2..4) {
    foreach ($col in @("E", "F")) {
        $overview.Range("$col$row").Value = $newText
    }
}

# zh-cn / de-de sheets: Status column is column C, rows 2-4
foreach ($ws in @($zhcn, $dede)) {
    foreach ($row in 2..4) {
        $ws.Range("C$row").Value = $newText
    }
}

# --- Narrow the now-shorter status columns -----------------------------
# The stored column width is driven by the pixel grid Excel snaps
# ColumnWidth onto; 12.5 characters lands on the same ~13.33-character
# rendered width that fitting "In Translation" produces (versus the wider
# ~17.2 characters that used to fit "Ready for handoff").
$newColumnWidth = 12.5

$overview.Range("E1:F1").ColumnWidth = $newColumnWidth
$zhcn.Range("C1").ColumnWidth = $newColumnWidth
$dede.Range("C1").ColumnWidth = $newColumnWidth
